$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F49").Value = 677.9061000678731
$ws.Range("G49").Value = 244.1388134236847
$ws.Range("I49").Value = 34
$ws.Range("K49").Value = 6.5

$ws.Range("F54").Value = 5773.99930088235
$ws.Range("G54").Value = 242.1331862562112
$ws.Range("I54").Value = 25
$ws.Range("K54").Value = 0.5

$ws.Range("F62").Value = 495.6342526451205
$ws.Range("H62").Value = 210.8411957790392

$ws.Range("F63").Value = 190
$ws.Range("G63").Value = 190
$ws.Range("H63").Value = 190

$ws.Range("F64").Value = 213.6597999528158
$ws.Range("H64").Value = 213.6597999528158

$ws.Range("F67").Value = 520.6087294613138
$ws.Range("H67").Value = 210.9419831262309

$ws.Range("F68").Value = 190
$ws.Range("G68").Value = 190
$ws.Range("H68").Value = 190

$ws.Range("F69").Value = 213.9446447551659
$ws.Range("H69").Value = 213.9446447551659

$ws.Range("F72").Value = 498.6846249564804
$ws.Range("H72").Value = 210.8547423725928

$ws.Range("F73").Value = 190
$ws.Range("G73").Value = 190
$ws.Range("H73").Value = 190

$ws.Range("F74").Value = 213.6824512306345
$ws.Range("H74").Value = 213.6824512306345

$ws.Range("F77").Value = 484.0617061442907
$ws.Range("H77").Value = 210.7861122794238

$ws.Range("F78").Value = 190
$ws.Range("G78").Value = 190
$ws.Range("H78").Value = 190

$ws.Range("F79").Value = 212.974716749215
$ws.Range("H79").Value = 212.974716749215

$ws.Range("F82").Value = 459.8310084052176
$ws.Range("H82").Value = 210.6477630245759

$ws.Range("F83").Value = 190
$ws.Range("G83").Value = 190
$ws.Range("H83").Value = 190

$ws.Range("F84").Value = 212.4819999512866
$ws.Range("H84").Value = 212.4819999512866

$ws.Range("F87").Value = 443.5375262346428
$ws.Range("H87").Value = 210.5316218693329

$ws.Range("F88").Value = 190
$ws.Range("G88").Value = 190
$ws.Range("H88").Value = 190

$ws.Range("F89").Value = 212.2633069458213
$ws.Range("H89").Value = 212.2633069458213

$ws.Range("F92").Value = 444.3408215816889
$ws.Range("H92").Value = 210.5378888793338

$ws.Range("F93").Value = 190
$ws.Range("G93").Value = 190
$ws.Range("H93").Value = 190

$ws.Range("F94").Value = 212.2143618186433
$ws.Range("H94").Value = 212.2143618186433

$ws.Range("F97").Value = 452.6956784699108
$ws.Range("H97").Value = 210.5996054923444

$ws.Range("F98").Value = 190
$ws.Range("G98").Value = 190
$ws.Range("H98").Value = 190

$ws.Range("F99").Value = 212.4980986425754
$ws.Range("H99").Value = 212.4980986425754

$ws.Range("F102").Value = 478.6261161211738
$ws.Range("H102").Value = 210.75801969462

$ws.Range("F103").Value = 190
$ws.Range("G103").Value = 190
$ws.Range("H103").Value = 190

$ws.Range("F104").Value = 212.9903954632742
$ws.Range("H104").Value = 212.9903954632742

$ws.Range("F107").Value = 500.5468872647287
$ws.Range("H107").Value = 210.8628280506697

$ws.Range("F108").Value = 190
$ws.Range("G108").Value = 190
$ws.Range("H108").Value = 190

$ws.Range("F109").Value = 213.2317286018183
$ws.Range("H109").Value = 213.2317286018183

$ws.Range("F112").Value = 498.0692570064664
$ws.Range("H112").Value = 210.8520401185984

$ws.Range("F113").Value = 190
$ws.Range("G113").Value = 190
$ws.Range("H113").Value = 190

$ws.Range("F114").Value = 213.2007658469566
$ws.Range("H114").Value = 213.2007658469566

$ws.Range("F117").Value = 496.846495637108
$ws.Range("H117").Value = 210.8466250273695

$ws.Range("F118").Value = 190
$ws.Range("G118").Value = 190
$ws.Range("H118").Value = 190

$ws.Range("F119").Value = 213.3627410968443
$ws.Range("H119").Value = 213.3627410968443

Write-Host "done"
